$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New unified comment text used for all data rows (column P)
$newComment = "calculated in EUCalc-Tool (Life Scenario), for details see script https://github.com/Sufficiency-Quantification/SufficiencyIndicators_EU-Calc/tree/main"

# Row 2
$ws.Range("J2").Value = "GJ/cap/year"
$ws.Range("L2").Value = 92.44023461698886
$ws.Range("P2").Value = $newComment

# Row 3
$ws.Range("J3").Value = "GJ/cap/year"
$ws.Range("L3").Value = 28.9900034845241
$ws.Range("P3").Value = $newComment

# Row 4
$ws.Range("J4").Value = "m²/cap"
$ws.Range("P4").Value = $newComment

# Row 5
$ws.Range("J5").Value = "m²/cap"
$ws.Range("L5").Value = 20.5748851717017
$ws.Range("P5").Value = $newComment

# Row 6
$ws.Range("J6").Value = "pkm/cap/year"
$ws.Range("P6").Value = $newComment

# Row 7
$ws.Range("J7").Value = "pkm/cap/year"
$ws.Range("P7").Value = $newComment

# Row 8
$ws.Range("J8").Value = "pkm/cap/year"
$ws.Range("P8").Value = $newComment

# Row 9
$ws.Range("J9").Value = "tkm/cap/year"
$ws.Range("L9").Value = 10597.54755692944
$ws.Range("P9").Value = $newComment

# Row 10
$ws.Range("J10").Value = "kcal meat/cap/day"
$ws.Range("P10").Value = $newComment

# Row 11
$ws.Range("J11").Value = "kcal/cap/year"
$ws.Range("P11").Value = $newComment

# Row 12
$ws.Range("J12").Value = "t/cap/year"
$ws.Range("L12").Value = 0.4091521796827101
$ws.Range("P12").Value = $newComment

# Row 13
$ws.Range("J13").Value = "t/cap/year"
$ws.Range("L13").Value = 0.5255786268492357
$ws.Range("P13").Value = $newComment

# Row 14
$ws.Range("J14").Value = "GJ/cap/year"
$ws.Range("L14").Value = 41.02228303672896
$ws.Range("P14").Value = $newComment

# Row 15
$ws.Range("J15").Value = "GJ/cap/year"
$ws.Range("L15").Value = 13.89764409309618
$ws.Range("P15").Value = $newComment

# Row 16
$ws.Range("J16").Value = "m²/cap"
$ws.Range("P16").Value = $newComment

# Row 17
$ws.Range("J17").Value = "m²/cap"
$ws.Range("L17").Value = 20.21879437211879
$ws.Range("P17").Value = $newComment

# Row 18
$ws.Range("J18").Value = "pkm/cap/year"
$ws.Range("P18").Value = $newComment

# Row 19
$ws.Range("J19").Value = "pkm/cap/year"
$ws.Range("P19").Value = $newComment

# Row 20
$ws.Range("J20").Value = "pkm/cap/year"
$ws.Range("P20").Value = $newComment

# Row 21
$ws.Range("J21").Value = "tkm/cap/year"
$ws.Range("L21").Value = 8123.025418079452
$ws.Range("P21").Value = $newComment

# Row 22
$ws.Range("J22").Value = "kcal meat/cap/day"
$ws.Range("P22").Value = $newComment

# Row 23
$ws.Range("J23").Value = "kcal/cap/year"
$ws.Range("P23").Value = $newComment

# Row 24
$ws.Range("J24").Value = "t/cap/year"
$ws.Range("L24").Value = 0.2050796055697152
$ws.Range("P24").Value = $newComment

# Row 25
$ws.Range("J25").Value = "t/cap/year"
$ws.Range("L25").Value = 0.314260784468168
$ws.Range("P25").Value = $newComment

# Row 26
$ws.Range("J26").Value = "GJ/cap/year"
$ws.Range("L26").Value = 74.36245531958068
$ws.Range("P26").Value = $newComment

# Row 27
$ws.Range("J27").Value = "GJ/cap/year"
$ws.Range("L27").Value = 20.72291162311338
$ws.Range("P27").Value = $newComment

# Row 28
$ws.Range("J28").Value = "m²/cap"
$ws.Range("P28").Value = $newComment

# Row 29
$ws.Range("J29").Value = "m²/cap"
$ws.Range("L29").Value = 14.69473041515453
$ws.Range("P29").Value = $newComment

# Row 30
$ws.Range("J30").Value = "pkm/cap/year"
$ws.Range("P30").Value = $newComment

# Row 31
$ws.Range("J31").Value = "pkm/cap/year"
$ws.Range("P31").Value = $newComment

# Row 32
$ws.Range("J32").Value = "pkm/cap/year"
$ws.Range("P32").Value = $newComment

# Row 33
$ws.Range("J33").Value = "tkm/cap/year"
$ws.Range("L33").Value = 7217.520129121297
$ws.Range("P33").Value = $newComment

# Row 34
$ws.Range("J34").Value = "kcal meat/cap/day"
$ws.Range("P34").Value = $newComment

# Row 35
$ws.Range("J35").Value = "kcal/cap/year"
$ws.Range("P35").Value = $newComment

# Row 36
$ws.Range("J36").Value = "t/cap/year"
$ws.Range("L36").Value = 0.3423666033624484
$ws.Range("P36").Value = $newComment

# Row 37
$ws.Range("J37").Value = "t/cap/year"
$ws.Range("L37").Value = 0.321444484124084
$ws.Range("P37").Value = $newComment

# Row 38
$ws.Range("J38").Value = "GJ/cap/year"
$ws.Range("L38").Value = 36.77920915105934
$ws.Range("P38").Value = $newComment

# Row 39
$ws.Range("J39").Value = "GJ/cap/year"
$ws.Range("L39").Value = 12.89167554005816
$ws.Range("P39").Value = $newComment

# Row 40
$ws.Range("J40").Value = "m²/cap"
$ws.Range("P40").Value = $newComment

# Row 41
$ws.Range("J41").Value = "m²/cap"
$ws.Range("L41").Value = 14.12437952123616
$ws.Range("P41").Value = $newComment

# Row 42
$ws.Range("J42").Value = "pkm/cap/year"
$ws.Range("P42").Value = $newComment

# Row 43
$ws.Range("J43").Value = "pkm/cap/year"
$ws.Range("P43").Value = $newComment

# Row 44
$ws.Range("J44").Value = "pkm/cap/year"
$ws.Range("P44").Value = $newComment

# Row 45
$ws.Range("J45").Value = "tkm/cap/year"
$ws.Range("L45").Value = 5411.159829960237
$ws.Range("P45").Value = $newComment

# Row 46
$ws.Range("J46").Value = "kcal meat/cap/day"
$ws.Range("P46").Value = $newComment

# Row 47
$ws.Range("J47").Value = "kcal/cap/year"
$ws.Range("P47").Value = $newComment

# Row 48
$ws.Range("J48").Value = "t/cap/year"
$ws.Range("L48").Value = 0.1856341912613852
$ws.Range("P48").Value = $newComment

# Row 49
$ws.Range("J49").Value = "t/cap/year"
$ws.Range("L49").Value = 0.2321482602218339
$ws.Range("P49").Value = $newComment
